# "update import and show detail"
#
# The import/header row on Sheet1 is reshuffled:
#   - B1 and C1 (nama_dengan_gelar / nama_tanpa_gelar) swap places
#   - E1's label "tempat,tanggal_lahir" is fixed to "tempat_tanggal_lahir"
#   - the active selection moves to B1 (detail view now opens on the
#     "show detail" column instead of A2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the B1 / C1 header labels.
$b1 = $ws.Range("B1").Value()
$c1 = $ws.Range("C1").Value()
$ws.Range("B1").Value = $c1
$ws.Range("C1").Value = $b1

# Correct the "tempat,tanggal_lahir" header text.
$ws.Range("E1").Value = "tempat_tanggal_lahir"

# Column C now holds the same-width label as column B; re-fit it so the
# two columns line up again (matches column B's best-fit width).
$ws.Columns.Item(3).ColumnWidth = 18.33

# Move / show the selection on the newly-relevant B1 cell.
$ws.Range("B1").Select()
